# Update cryptocurrency price/volume figures (Wed Jan 24 19:55:47 UTC 2024 data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings; some values (single-dot decimals) would be
# auto-coerced to numbers by a plain Value assignment, silently retyping the
# cell. Force the cell to Text first, write the literal string, then strip the
# Text number-format back off (ClearFormats) so the cell style/format matches
# the original (unstyled) cell exactly while the stored value stays a string.
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue 'D2' '39.700.09'
$ws.Range('E2').Value = '  +1.09%  '
Set-TextValue 'D3' '2.210.32'
$ws.Range('E3').Value = '  +0.87%  '
$ws.Range('E4').Value = '  -0.03%  '
Set-TextValue 'D5' '292.21'
$ws.Range('E5').Value = '  -1.09%  '
Set-TextValue 'D6' '86.07'
$ws.Range('E6').Value = '  +5.92%  '
$ws.Range('E7').Value = '  +1.05%  '
$ws.Range('E8').Value = '  -0.11%  '
$ws.Range('E9').Value = '  +1.43%  '
Set-TextValue 'D10' '30.56'
$ws.Range('E10').Value = '  +5.38%  '
$ws.Range('E11').Value = '  +2.00%  '
Set-TextValue 'D12' '47.52'
$ws.Range('E12').Value = '  +1.04%  '
$ws.Range('E13').Value = '  +1.40%  '
$ws.Range('E14').Value = '  +1.75%  '
Set-TextValue 'D15' '2.553.17'
$ws.Range('E15').Value = '  +0.99%  '
$ws.Range('E16').Value = '  +0.58%  '
Set-TextValue 'D17' '2.217.30'
$ws.Range('E17').Value = '  +1.63%  '
$ws.Range('E18').Value = '  +3.09%  '
Set-TextValue 'D19' '39.654.95'
$ws.Range('E19').Value = '  +1.29%  '
$ws.Range('E20').Value = '  +1.40%  '
$ws.Range('E21').Value = '  +9.56%  '
Set-TextValue 'D22' '5.79'
$ws.Range('E22').Value = '  +1.83%  '
Set-TextValue 'D23' '65.56'
$ws.Range('E23').Value = '  +1.15%  '
Set-TextValue 'D24' '235.31'
$ws.Range('E24').Value = '  +4.55%  '
$ws.Range('E25').Value = '  -0.19%  '
$ws.Range('E26').Value = '  +2.72%  '
$ws.Range('E27').Value = '  +2.41%  '
Set-TextValue 'D28' '22.70'
$ws.Range('E28').Value = '  +1.04%  '
Set-TextValue 'D29' '2.11'
$ws.Range('E29').Value = '  -2.49%  '
Set-TextValue 'D30' '9.23'
$ws.Range('E30').Value = '  +2.08%  '
Set-TextValue 'D31' '32.75'
$ws.Range('E31').Value = '  +3.24%  '
Set-TextValue 'D32' '151.37'
$ws.Range('E32').Value = '  +1.00%  '
Set-TextValue 'D33' '0.998'
$ws.Range('E33').Value = '  -0.18%  '
$ws.Range('E34').Value = '  +2.72%  '
$ws.Range('E35').Value = '  +3.66%  '
$ws.Range('E36').Value = '  +2.12%  '
$ws.Range('E37').Value = '  +6.72%  '
$ws.Range('E38').Value = '  +1.72%  '
Set-TextValue 'D39' '15.83'
$ws.Range('E39').Value = '  +3.53%  '
Set-TextValue 'D40' '0.0985'
$ws.Range('E40').Value = '  +2.56%  '
$ws.Range('E41').Value = '  +4.02%  '
Set-TextValue 'D42' '2.067.41'
$ws.Range('E42').Value = '  +9.30%  '
$ws.Range('E43').Value = '  +5.21%  '
$ws.Range('E44').Value = '  +3.22%  '
Set-TextValue 'D45' '9.92'
$ws.Range('E45').Value = '  +11.24%  '
$ws.Range('E46').Value = '  +0.35%  '
Set-TextValue 'D47' '17.72'
$ws.Range('E47').Value = '  +10.44%  '
$ws.Range('E48').Value = '  +0.22%  '
Set-TextValue 'D49' '2.430.20'
$ws.Range('E49').Value = '  +1.32%  '
Set-TextValue 'D50' '71.10'
$ws.Range('E50').Value = '  +0.05%  '
Set-TextValue 'D51' '88.91'
$ws.Range('E51').Value = '  +2.69%  '
